$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Change 1 ---------------------------------------------------------
# Row 8 (Vartotojas iveda filmo pavadinima ...), column 2 is an empty
# "Sistemos reakcija" paragraph. Fill it in with "1.1 Sistema tikrina"
# (split across a few runs, matching the original authoring), the
# relocated _GoBack bookmark, and the trailing " ivestus duomenis."
$cell1 = $t.Cell(8, 2)
$p1 = $cell1.Range.Paragraphs.Item(1)
$r1 = $p1.Range
$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00DA3747" w:rsidRDefault="00DA3747" w:rsidP="003F32E2"><w:pPr><w:rPr><w:lang w:val="lt-LT"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="lt-LT"/></w:rPr><w:t>1</w:t></w:r><w:r><w:rPr><w:lang w:val="lt-LT"/></w:rPr><w:t xml:space="preserve">.1 </w:t></w:r><w:r><w:rPr><w:lang w:val="lt-LT"/></w:rPr><w:t>Sistema tikrina</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:lang w:val="lt-LT"/></w:rPr><w:t xml:space="preserve"> įvestus duomenis.</w:t></w:r></w:p>'
$r1.InsertXML($xml1)

# --- Change 2 -----------------------------------------------------------
# Row 9 ("Vartotojas nuspaudzia filmo kureju pridejimo mygtuka."), column
# 2 currently reads "2.1 Sistema pateikia filmo kūrėjų paiešką" with no
# trailing period. Add one, as its own run.
$cell2 = $t.Cell(9, 2)
$p2 = $cell2.Range.Paragraphs.Item(1)
$r2 = $p2.Range
$xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00DA3747" w:rsidRDefault="00674ED4" w:rsidP="003F32E2"><w:pPr><w:rPr><w:lang w:val="lt-LT"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="lt-LT"/></w:rPr><w:t>2.1 Sistema pateikia filmo kūrėjų paiešką</w:t></w:r><w:r><w:rPr><w:lang w:val="lt-LT"/></w:rPr><w:t>.</w:t></w:r></w:p>'
$r2.InsertXML($xml2)

# --- Change 3 -------------------------------------------------------------
# Row 11 ("Nuspaudziamas filmo pridejimo mygtukas"), column 2 ends with
# "4.1 Sistema issiuncia ... patvirtinti pridejima." followed by the
# _GoBack bookmark - that bookmark moved up to Change 1, so drop it here.
$cell3 = $t.Cell(11, 2)
$p3 = $cell3.Range.Paragraphs.Item(1)
$r3 = $p3.Range
$xml3 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="009119C9" w:rsidRDefault="009119C9" w:rsidP="003F32E2"><w:pPr><w:rPr><w:lang w:val="lt-LT"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="lt-LT"/></w:rPr><w:t>4.1 Sistema išsiunčia pranešimus filmų kūrėjams, kurie buvo pridėti prie filmo, patvirtinti pridėjimą.</w:t></w:r></w:p>'
$r3.InsertXML($xml3)
